$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: give B2 the bordered/shaded style and add matching C2/D2 ---
$ws.Range("B3").Copy()
$ws.Range("B2:D2").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B2:D2").Value = 0

# --- Update existing Well-1 survey rows (3 & 4) ---
$ws.Range("C3").Value = 9
$ws.Range("D3").Value = 35
$ws.Range("C4").Value = 50
$ws.Range("D4").Value = 37

# --- New Well-2 block (rows 5-7) ---
$ws.Range("A5").Value = "Well-2"
$ws.Range("A6").Value = "Well-2"
$ws.Range("A7").Value = "Well-2"

# Row 5 mirrors row 2's "zeroed" style (border + shading)
$ws.Range("B3").Copy()
$ws.Range("B5:D5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B5:D5").Value = 0

# Rows 6-7 use the shaded style with no border (new style)
$ws.Range("B3").Copy()
$ws.Range("B6:D7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B6:D7").Borders.LineStyle = -4142  # xlLineStyleNone

$ws.Range("B6").Value = 800
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 20
$ws.Range("B7").Value = 985
$ws.Range("C7").Value = 30
$ws.Range("D7").Value = 25

$excel.CutCopyMode = $false

# --- View state ---
$ws.Application.ActiveWindow.Zoom = 70
[void]$ws.Range("D9").Select()
